$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the top
#    of the data, shift every other quarter down by one row, and add a new
#    last row for 2020-Q4.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$totalRows = @(
    @(0, "2022-Q4", 9, 0.99),
    @(1, "2022-Q3", 9, 1.05),
    @(2, "2022-Q2", 6, 1.15),
    @(3, "2022-Q1", 7, 1.27),
    @(4, "2021-Q4", 9, 2.8),
    @(5, "2021-Q3", 30, 11.51),
    @(6, "2021-Q2", 14, 1.49),
    @(7, "2021-Q1", 5, 0.4),
    @(8, "2020-Q4", 2, 0.1)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

# Row 10 (2020-Q4) is brand new - copy column-A's formatting (bold + border)
# from the previous last row (row 9) so it matches the rest of the column.
$fmtSrc = $total.Cells.Item(9, 1)
$fmtSrc.Copy()
$total.Cells.Item(10, 1).PasteSpecial(-4122)

Write-Output "updated total sheet"

# ---------------------------------------------------------------------------
# 2. Add the new "2022-Q4" fund-holdings sheet. The easiest way to get an
#    exact match of formatting (fonts/borders/column widths/etc.) is to
#    duplicate the existing "2022-Q3" sheet (which carries that formatting)
#    and then overwrite the numbers that changed for the new quarter. The
#    duplicate is placed right after "总计", matching the target tab order.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Fund code (B) / fund name (C) / row index (A) stay the same as 2022-Q3;
# only the size/position/weight/value/rank metrics (D-H) change.
# D, E, F, G were authored as text in the source workbook (pandas export),
# so we re-enter them with a leading apostrophe to keep them as text
# instead of letting Excel auto-convert them to numbers.
$q4Rows = @(
    @(2,  "13.53", "99.48", "2.67", "0.3613", 9),
    @(3,  "8.97",  "94.55", "2.64", "0.2368", 8),
    @(4,  "10.97", "94.25", "1.54", "0.1689", 7),
    @(5,  "3.28",  "94.55", "2.64", "0.0866", 8),
    @(6,  "3.18",  "92.38", "2.57", "0.0817", 8),
    @(7,  "3.29",  "94.25", "1.54", "0.0507", 7),
    @(8,  "0.08",  "91.91", "2.49", "0.0020", 10),
    @(9,  "0.02",  "91.91", "2.49", "0.0005", 10),
    @(10, "0.02",  "92.38", "2.57", "0.0005", 8)
)

foreach ($row in $q4Rows) {
    $r = $row[0]
    $q4.Cells.Item($r, 4).Value = "'" + $row[1]
    $q4.Cells.Item($r, 5).Value = "'" + $row[2]
    $q4.Cells.Item($r, 6).Value = "'" + $row[3]
    $q4.Cells.Item($r, 7).Value = "'" + $row[4]
    $q4.Cells.Item($r, 8).Value = $row[5]
}

Write-Output "added 2022-Q4 sheet"
